# Sample Project / Main.xlsx — "Rules" sheet
# B11 previously held the text "R40"; it must now hold the text "1"
# (still a plain text value, not a number), with its original cell
# formatting (style) left completely untouched.
#
# A bare assignment like $cell.Value = "1" would be auto-recognised by
# Excel as a number. To keep it textual we have to go through Excel's
# normal "quote prefix" mechanism (the same thing that happens when a
# user types '1 into a cell) - but that by itself also stamps a fresh
# "quote prefixed" style onto the cell. So: stash the cell's current
# formatting on a scratch cell, make the text assignment, then paste the
# original formatting back on top so the cell's style is restored
# exactly as it was. Finally remove the scratch column so the sheet's
# used range/dimensions are left exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

$scratchCol = 26                 # column Z - far outside the sheet's real data
$scratch = $ws.Cells.Item(1, $scratchCol)

# 1) Remember the cell's current formatting on the scratch cell.
$scratch.Value = "x"
$cell.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 2) Assign the new text value, forcing text (not numeric) interpretation.
$cell.Value = "'1"

# 3) Restore the original formatting (clears the quote-prefix style bump).
$scratch.Copy() | Out-Null
$cell.PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# 4) Clean up the scratch column so no stray cells/dimensions remain.
$ws.Columns.Item($scratchCol).Delete()
